$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Fri Sep 29 11:37:52 EDT 2023"
$ws.Range("B3").Value = "Fri Sep 29 11:38:04 EDT 2023"
$ws.Range("B4").Value = "Fri Sep 29 11:38:17 EDT 2023"
$ws.Range("B5").Value = "Fri Sep 29 11:38:29 EDT 2023"
